$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
# The "Periodo Mora" column (E16:E22) is re-keyed to the ascending run of
# periods 2205-2211 (part 1 of the new account-statement data), and the
# "Valor Mora" figures for the 2205/2211 periods (rows 16 and 22) are
# updated to reflect the refreshed database values.
$ws.Range("E16").Value = "2205"
$ws.Range("E17").Value = "2206"
$ws.Range("E18").Value = "2207"
$ws.Range("E19").Value = "2208"
$ws.Range("E20").Value = "2209"
$ws.Range("E21").Value = "2210"
$ws.Range("E22").Value = "2211"

$ws.Range("F16").Value = 40000
$ws.Range("F22").Value = 28000
